# Autogenerated on Sun Feb 01 2015 22:24:41 GMT-0500 (Eastern Standard Time)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from "Data" to "Summary"
$ws.Name = "Summary"

# Add a new named cell style "title_" (bold + underlined), mirroring the
# existing "title" style's font size/face, appended to the style gallery.
$newStyle = $wb.Styles.Add("title_")
$newStyle.Font.Bold = $true
$newStyle.Font.Underline = $true

# Remove row 5 (the Micro / SMEs / MSMEs header row) entirely.
$ws.Rows.Item(5).Delete()
